# Append 5 new daily-report rows (132-136) to the Arequipa COVID data sheet,
# matching the style pattern already used for the last few rows of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row data: date-serial, total_muestras, casos_positivos, casos_negativos,
# espera_resultado, defunciones, recuperados, hospitalizados_positivos,
# hospitalizados_sospechosos, uci_positivos, uci_sospechosos, trauma_shock,
# defunciones_minsa, defunciones_essalud, defunciones_clinicas,
# defunciones_domicilio, defunciones_sanidades
# ---------------------------------------------------------------------------
$rows = @(
    @(132, 44082, 661464, 109736, 551664, 64,  1821, 86502, 402, 96,  30, 0, 10, 693, 1039, 38, 41,  10),
    @(133, 44083, 668597, 111085, 557470, 42,  1834, 88366, 355, 126, 31, 0, 10, 695, 1050, 38, 41,  10),
    @(134, 44084, 676147, 112487, 563600, 60,  1844, 89773, 381, 127, 30, 0, 10, 696, 1059, 38, 41,  10),
    @(135, 44085, 683346, 113577, 569612, 157, 1862, 90618, 604, 132, 72, 0, 10, 706, 1067, 38, 41,  10),
    @(136, 44086, 689465, 114622, 574668, 175, 1948, 90394, 604, 132, 72, 0, 10, 711, 1073, 38, 116, 10)
)

# Write all of the raw values first, one cell at a time (array/range bulk
# assignment isn't reliable in this host, but per-cell Value writes are).
foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 17; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}

# ---------------------------------------------------------------------------
# Formatting: seed every new row with the formatting of the row above it
# (row 131), which carries the "date" format in column A (yyyy-mm-dd) and
# the plain numeric format across columns B:Q -- then layer on the
# right-alignment that was applied to the newly-added rows.
# ---------------------------------------------------------------------------
$ws.Range("A131:Q131").Copy() | Out-Null
for ($r = 132; $r -le 136; $r++) {
    $ws.Range("A" + $r + ":Q" + $r).PasteSpecial(-4122) | Out-Null
}

# Column A of every new row is right-aligned (date column).
for ($r = 132; $r -le 136; $r++) {
    $ws.Range("A" + $r).HorizontalAlignment = -4152
}

# Row 132 additionally has every other column (B:Q) right-aligned.
$ws.Range("B132:Q132").HorizontalAlignment = -4152

$excel.CutCopyMode = $false

Write-Host "Added rows 132-136"
